$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember old B9 text ("TIPO") before we shift columns.
$oldB9 = $ws.Range("B9").Value2

# Insert a new column "inside" the B:E (etc) merged ranges so that Excel
# naturally expands those merges to cover the new column, matching how
# the workbook was actually edited (insert column, then fill in the new
# "ESTABLECIMIENTO" picker column before the old column B).
$ws.Columns("C").EntireColumn.Insert()

# --- Row 4 (EMPRESA / ESTABLECIMEINTO labels) -------------------------
# A4 used to read "ESTABLECIMEINTO :" -- it now becomes the new "EMPRESA:"
# label (typed first), and the old "ESTABLECIMEINTO :" label moves over to
# G4 (typed second) -- this ordering matches the shared-string indices
# created by the original edit.
$ws.Range("A4").Value = "EMPRESA:"
$ws.Range("G4").Value = "ESTABLECIMEINTO :"

# Copy A4's label style (bold / filled / bordered) onto G4.
$ws.Range("A4").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New merged "value" area next to the relocated ESTABLECIMEINTO label.
$ws.Range("H4:K4").Merge()

# --- Row 9 (header row, no merges) -----------------------------------
# The insert above left the new blank cell at C9 (pushing old C9.. to D9..),
# but logically the new column belongs at B9 with the old "TIPO" header
# shifting right into C9.
$ws.Range("C9").Value = $oldB9
$ws.Range("B9").Value = "ESTABLECIMIENTO"

# Column B (the new ESTABLECIMIENTO picker column) needs its own width.
$ws.Columns("B").ColumnWidth = 15.999

# Selection ends up on the newly typed header cell, like in the source edit.
$ws.Range("B9").Select()

Write-Host "edit applied"
